# Removes the Aspose.Slides "Evaluation only" watermark textboxes that were
# left over on the two slides, and refreshes the cached date shown by the
# date placeholder on the slide layout (recorded the last time the deck was
# saved).

$p = $ppt.ActivePresentation

# --- Remove the "Evaluation only." watermark textbox from every slide ----
for ($si = 1; $si -le $p.Slides.Count; $si++) {
    $slide = $p.Slides.Item($si)
    for ($i = $slide.Shapes.Count; $i -ge 1; $i--) {
        $shp = $slide.Shapes.Item($i)
        if ($shp.Name -eq "TextBox") {
            $shp.Delete()
        }
    }
}

# --- Refresh the cached "datetimeFigureOut" text on the title layout -----
$layout = $p.Slides.Item(1).CustomLayout
for ($i = 1; $i -le $layout.Shapes.Count; $i++) {
    $phShape = $layout.Shapes.Item($i)
    if ($phShape.Name -eq "Date Placeholder 1") {
        $phShape.TextFrame.TextRange.Text = "5/26/2016"
    }
}
